# Generate Report for Handback
# Refresh the handback-status report's timestamps: the "Latest HO Xliff
# Generate Date" on the Overview sheet, and the per-locale "Correspond
# Handoff/Handback Datetime" columns on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date (row 2 = the 41b5d314-...md file)
$wsOverview.Range("G2").Value = "2016-10-18 12:52:02"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-10-18 12:51:50"
$wsZhCn.Range("K2").Value = "2016-10-18 12:52:29"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-10-18 12:52:02"
$wsDeDe.Range("K2").Value = "2016-10-18 12:52:46"
